# Scheduled market-price refresh for the Kujata_Profits leve-profit tables.
# Updates currentAveragePrice* / LeveProfit* columns (H, I, J, K, L, M, N) for
# the rows whose underlying item prices changed, across all eight crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5599.4287
$ws.Range("I76").Value = 4732.6665
$ws.Range("K76").Value = 4732.6665
$ws.Range("M76").Value = -4417.6665

$ws.Range("H79").Value = 5599.4287
$ws.Range("I79").Value = 4732.6665
$ws.Range("K79").Value = 4732.6665
$ws.Range("M79").Value = -3640.6665

$ws.Range("H107").Value = 4218.55
$ws.Range("I107").Value = 3317.0625
$ws.Range("K107").Value = 3317.0625
$ws.Range("M107").Value = -1397.0625

$ws.Range("H137").Value = 2041.8363
$ws.Range("I137").Value = 1656.7778
$ws.Range("J137").Value = 2413.1428
$ws.Range("K137").Value = 4970.3334
$ws.Range("L137").Value = 7239.428400000001
$ws.Range("M137").Value = -2420.3334
$ws.Range("N137").Value = -12339.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10266.237
$ws.Range("I32").Value = 7604.855
$ws.Range("K32").Value = 7604.855
$ws.Range("M32").Value = -7317.855

$ws.Range("H61").Value = 90910690
$ws.Range("I61").Value = 125001190
$ws.Range("K61").Value = 125001190
$ws.Range("M61").Value = -125000978

$ws.Range("H63").Value = 33335564
$ws.Range("I63").Value = 2282
$ws.Range("K63").Value = 2282
$ws.Range("M63").Value = -1596

$ws.Range("H66").Value = 33335564
$ws.Range("I66").Value = 2282
$ws.Range("K66").Value = 11410
$ws.Range("M66").Value = -7978

$ws.Range("H74").Value = 1442.9474
$ws.Range("I74").Value = 977.41174
$ws.Range("K74").Value = 977.41174
$ws.Range("M74").Value = -103.41174

$ws.Range("H77").Value = 1442.9474
$ws.Range("I77").Value = 977.41174
$ws.Range("K77").Value = 4887.0587
$ws.Range("M77").Value = -519.0586999999996

$ws.Range("H88").Value = 2255.4443
$ws.Range("J88").Value = 2114.2856
$ws.Range("L88").Value = 2114.2856
$ws.Range("N88").Value = -2926.2856

$ws.Range("H91").Value = 2255.4443
$ws.Range("J91").Value = 2114.2856
$ws.Range("L91").Value = 2114.2856
$ws.Range("N91").Value = -4922.2856

$ws.Range("H136").Value = 90910690
$ws.Range("I136").Value = 125001190
$ws.Range("K136").Value = 375003570
$ws.Range("M136").Value = -375001020

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1391.2
$ws.Range("I107").Value = 1149.875
$ws.Range("K107").Value = 1149.875
$ws.Range("M107").Value = 770.125

$ws.Range("H134").Value = 3616.389
$ws.Range("I134").Value = 829.14703
$ws.Range("K134").Value = 2487.44109
$ws.Range("M134").Value = 47.5589100000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 76924410
$ws.Range("I16").Value = 111112550
$ws.Range("J16").Value = 1081.5
$ws.Range("K16").Value = 111112550
$ws.Range("L16").Value = 1081.5
$ws.Range("M16").Value = -111112263
$ws.Range("N16").Value = -1655.5

$ws.Range("H31").Value = 1264.4038
$ws.Range("I31").Value = 1249.7059
$ws.Range("K31").Value = 1249.7059
$ws.Range("M31").Value = -954.7058999999999

$ws.Range("H34").Value = 1264.4038
$ws.Range("I34").Value = 1249.7059
$ws.Range("K34").Value = 1249.7059
$ws.Range("M34").Value = -1047.7059

$ws.Range("H53").Value = 29000
$ws.Range("J53").Value = 29000
$ws.Range("L53").Value = 29000
$ws.Range("N53").Value = -30214

$ws.Range("H105").Value = 756.44446
$ws.Range("I105").Value = 726
$ws.Range("K105").Value = 726
$ws.Range("M105").Value = 1021

$ws.Range("H107").Value = 609.25
$ws.Range("I107").Value = 401.22223
$ws.Range("J107").Value = 1233.3334
$ws.Range("K107").Value = 401.22223
$ws.Range("L107").Value = 1233.3334
$ws.Range("M107").Value = 1518.77777
$ws.Range("N107").Value = -5073.3334

$ws.Range("H112").Value = 38450.332
$ws.Range("J112").Value = 38450.332
$ws.Range("L112").Value = 38450.332
$ws.Range("N112").Value = -41404.332

$ws.Range("H113").Value = 76924410
$ws.Range("I113").Value = 111112550
$ws.Range("J113").Value = 1081.5
$ws.Range("K113").Value = 111112550
$ws.Range("L113").Value = 1081.5
$ws.Range("M113").Value = -111110380
$ws.Range("N113").Value = -5421.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 512.08826
$ws.Range("J5").Value = 1287
$ws.Range("L5").Value = 3861
$ws.Range("N5").Value = -4085

$ws.Range("H107").Value = 7191.8667
$ws.Range("J107").Value = 9636.182000000001
$ws.Range("L107").Value = 28908.546
$ws.Range("N107").Value = -32748.546

$ws.Range("H131").Value = 34535524
$ws.Range("J131").Value = 84605.55499999999
$ws.Range("L131").Value = 253816.665
$ws.Range("N131").Value = -263896.665

$ws.Range("H135").Value = 512.08826
$ws.Range("J135").Value = 1287
$ws.Range("L135").Value = 11583
$ws.Range("N135").Value = -16653

$ws.Range("H139").Value = 2361.7083
$ws.Range("I139").Value = 2332.5217
$ws.Range("J139").Value = 3033
$ws.Range("K139").Value = 6997.5651
$ws.Range("L139").Value = 9099
$ws.Range("M139").Value = -1857.5651
$ws.Range("N139").Value = -19379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 34619020
$ws.Range("I70").Value = 31253714
$ws.Range("K70").Value = 31253714
$ws.Range("M70").Value = -31253444

$ws.Range("H73").Value = 34619020
$ws.Range("I73").Value = 31253714
$ws.Range("K73").Value = 31253714
$ws.Range("M73").Value = -31252778

$ws.Range("H80").Value = 4557.143
$ws.Range("J80").Value = 5016.6665
$ws.Range("L80").Value = 5016.6665
$ws.Range("N80").Value = -7012.6665

$ws.Range("H83").Value = 4557.143
$ws.Range("J83").Value = 5016.6665
$ws.Range("L83").Value = 25083.3325
$ws.Range("N83").Value = -35067.3325

$ws.Range("H102").Value = 2482.7932
$ws.Range("I102").Value = 1933.3125
$ws.Range("J102").Value = 3159.077
$ws.Range("K102").Value = 1933.3125
$ws.Range("L102").Value = 3159.077
$ws.Range("M102").Value = -311.3125
$ws.Range("N102").Value = -6403.077

$ws.Range("H122").Value = 1098
$ws.Range("I122").Value = 1098
$ws.Range("K122").Value = 3294
$ws.Range("M122").Value = -844

$ws.Range("H130").Value = 39280
$ws.Range("J130").Value = 39280
$ws.Range("L130").Value = 39280
$ws.Range("N130").Value = -49320

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5467.7144
$ws.Range("I40").Value = 3097.25
$ws.Range("J40").Value = 8628.333000000001
$ws.Range("K40").Value = 3097.25
$ws.Range("L40").Value = 8628.333000000001
$ws.Range("M40").Value = -2961.25
$ws.Range("N40").Value = -8900.333000000001

$ws.Range("H74").Value = 21098.5
$ws.Range("I74").Value = 13197
$ws.Range("J74").Value = 29000
$ws.Range("K74").Value = 13197
$ws.Range("L74").Value = 29000
$ws.Range("M74").Value = -12199
$ws.Range("N74").Value = -30996

$ws.Range("H77").Value = 21098.5
$ws.Range("I77").Value = 13197
$ws.Range("J77").Value = 29000
$ws.Range("K77").Value = 39591
$ws.Range("L77").Value = 87000
$ws.Range("M77").Value = -34599
$ws.Range("N77").Value = -96984

$ws.Range("H136").Value = 1831.5
$ws.Range("I136").Value = 1077
$ws.Range("J136").Value = 2586
$ws.Range("K136").Value = 3231
$ws.Range("L136").Value = 7758
$ws.Range("M136").Value = -681
$ws.Range("N136").Value = -12858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 11075
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H107").Value = 512
$ws.Range("I107").Value = 449.5
$ws.Range("J107").Value = 699.5
$ws.Range("K107").Value = 1348.5
$ws.Range("L107").Value = 2098.5
$ws.Range("M107").Value = 571.5
$ws.Range("N107").Value = -5938.5

$ws.Range("H132").Value = 2504.7354
$ws.Range("J132").Value = 4399.6
$ws.Range("L132").Value = 13198.8
$ws.Range("N132").Value = -18258.8
